$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 195
$ws.Range("I18").Value = 195
$ws.Range("K18").Value = 195
$ws.Range("M18").Value = 89
$ws.Range("H88").Value = 2359.6
$ws.Range("I88").Value = 2183
$ws.Range("J88").Value = 2624.5
$ws.Range("K88").Value = 2183
$ws.Range("L88").Value = 2624.5
$ws.Range("M88").Value = -1777
$ws.Range("N88").Value = -3436.5
$ws.Range("H91").Value = 2359.6
$ws.Range("I91").Value = 2183
$ws.Range("J91").Value = 2624.5
$ws.Range("K91").Value = 2183
$ws.Range("L91").Value = 2624.5
$ws.Range("M91").Value = -779
$ws.Range("N91").Value = -5432.5
$ws.Range("H112").Value = 1815.4286
$ws.Range("J112").Value = 1991.6
$ws.Range("L112").Value = 5974.799999999999
$ws.Range("N112").Value = -8190.799999999999
$ws.Range("H134").Value = 124900
$ws.Range("J134").Value = 124900
$ws.Range("L134").Value = 124900
$ws.Range("N134").Value = -135040
$ws.Range("H138").Value = 4497.222
$ws.Range("I138").Value = 3991.1538
$ws.Range("J138").Value = 4783.2607
$ws.Range("K138").Value = 11973.4614
$ws.Range("L138").Value = 14349.7821
$ws.Range("M138").Value = -6833.4614
$ws.Range("N138").Value = -24629.7821
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11421.625
$ws.Range("I32").Value = 11185.096
$ws.Range("K32").Value = 11185.096
$ws.Range("M32").Value = -10898.096
$ws.Range("H37").Value = 37998.75
$ws.Range("J37").Value = 37998.75
$ws.Range("L37").Value = 37998.75
$ws.Range("N37").Value = -38544.75
$ws.Range("H44").Value = 39995
$ws.Range("J44").Value = 39995
$ws.Range("L44").Value = 39995
$ws.Range("N44").Value = -40971
$ws.Range("H55").Value = 58053
$ws.Range("J55").Value = 58053
$ws.Range("L55").Value = 58053
$ws.Range("N55").Value = -58683
$ws.Range("H138").Value = 99894.5
$ws.Range("J138").Value = 99894.5
$ws.Range("L138").Value = 99894.5
$ws.Range("N138").Value = -110174.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1116.9166
$ws.Range("I107").Value = 940.5
$ws.Range("K107").Value = 940.5
$ws.Range("M107").Value = 979.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -350
$ws.Range("H31").Value = 2371
$ws.Range("I31").Value = 2314
$ws.Range("K31").Value = 2314
$ws.Range("M31").Value = -2019
$ws.Range("H34").Value = 2371
$ws.Range("I34").Value = 2314
$ws.Range("K34").Value = 2314
$ws.Range("M34").Value = -2112
$ws.Range("H62").Value = 4371
$ws.Range("J62").Value = 4144
$ws.Range("L62").Value = 4144
$ws.Range("N62").Value = -5392
$ws.Range("H65").Value = 4371
$ws.Range("J65").Value = 4144
$ws.Range("L65").Value = 20720
$ws.Range("N65").Value = -26960
$ws.Range("H107").Value = 2258.2
$ws.Range("I107").Value = 1572.75
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1572.75
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 347.25
$ws.Range("N107").Value = -8840
$ws.Range("H132").Value = 4141.2354
$ws.Range("I132").Value = 4026.8667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 12080.6001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -9550.6001
$ws.Range("N132").Value = -20057
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16446404
$ws.Range("I4").Value = 16751407
$ws.Range("J4").Value = 16080400
$ws.Range("K4").Value = 50254221
$ws.Range("L4").Value = 48241200
$ws.Range("M4").Value = -50254109
$ws.Range("N4").Value = -48241424
$ws.Range("H23").Value = 357.25
$ws.Range("J23").Value = 357.25
$ws.Range("L23").Value = 1071.75
$ws.Range("N23").Value = -1541.75
$ws.Range("H113").Value = 3497.625
$ws.Range("J113").Value = 3426.5715
$ws.Range("L113").Value = 10279.7145
$ws.Range("N113").Value = -14619.7145
$ws.Range("H121").Value = 2982.8333
$ws.Range("I121").Value = 560
$ws.Range("K121").Value = 1680
$ws.Range("M121").Value = -370
$ws.Range("H131").Value = 1915.375
$ws.Range("I131").Value = 1448
$ws.Range("K131").Value = 4344
$ws.Range("M131").Value = 696
$ws.Range("H140").Value = 1528
$ws.Range("I140").Value = 1528
$ws.Range("K140").Value = 4584
$ws.Range("M140").Value = 596
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.28570999999999
$ws.Range("J2").Value = 50
$ws.Range("L2").Value = 50
$ws.Range("N2").Value = -276
$ws.Range("H102").Value = 2956
$ws.Range("I102").Value = 2938.6
$ws.Range("K102").Value = 2938.6
$ws.Range("M102").Value = -1316.6
$ws.Range("H132").Value = 4561.8945
$ws.Range("I132").Value = 4581.75
$ws.Range("J132").Value = 4527.857
$ws.Range("K132").Value = 13745.25
$ws.Range("L132").Value = 13583.571
$ws.Range("M132").Value = -11215.25
$ws.Range("N132").Value = -18643.571
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6087.778
$ws.Range("I22").Value = 5631
$ws.Range("K22").Value = 5631
$ws.Range("M22").Value = -5336
$ws.Range("H27").Value = 6087.778
$ws.Range("I27").Value = 5631
$ws.Range("K27").Value = 5631
$ws.Range("M27").Value = -5524
$ws.Range("H55").Value = 399.55554
$ws.Range("I55").Value = 404.2
$ws.Range("K55").Value = 404.2
$ws.Range("M55").Value = -231.2
$ws.Range("H100").Value = 1599.9
$ws.Range("I100").Value = 1078.5714
$ws.Range("J100").Value = 2816.3333
$ws.Range("K100").Value = 1078.5714
$ws.Range("L100").Value = 2816.3333
$ws.Range("M100").Value = -537.5714
$ws.Range("N100").Value = -3898.3333
$ws.Range("H136").Value = 5576.2856
$ws.Range("I136").Value = 5907.4
$ws.Range("J136").Value = 4748.5
$ws.Range("K136").Value = 17722.2
$ws.Range("L136").Value = 14245.5
$ws.Range("M136").Value = -15172.2
$ws.Range("N136").Value = -19345.5
